# Rename the inline Pearson/BTec logo pictures that live in the document's
# header/footer stories.
#   - footer1 / footer2 : Pearson logo  image2.png -> image1.png
#   - header2           : BTec logo     image1.jpg -> image2.jpg
#
# These pictures are not part of the main body, so they are reached through
# Sections(1).Headers/Footers rather than $d.InlineShapes.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footers: Pearson logo, image2.png -> image1.png ---------------------
for ($fi = 1; $fi -le $sec.Footers.Count; $fi++) {
    $ftr = $sec.Footers.Item($fi)
    if ($ftr.Exists) {
        for ($si = 1; $si -le $ftr.Range.InlineShapes.Count; $si++) {
            $shp = $ftr.Range.InlineShapes.Item($si)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}

# --- Headers: BTec logo, image1.jpg -> image2.jpg -------------------------
for ($hi = 1; $hi -le $sec.Headers.Count; $hi++) {
    $hdr = $sec.Headers.Item($hi)
    if ($hdr.Exists) {
        for ($si = 1; $si -le $hdr.Range.InlineShapes.Count; $si++) {
            $shp = $hdr.Range.InlineShapes.Item($si)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image2.jpg"
            }
        }
    }
}
